# "more small fixes to sorting intro slides"
#
# Applies the three text corrections from the commit that are reachable
# through the PowerPoint object model:
#   - Slide 12: "A[0..j-1]" -> "A[1..j-1]"
#   - Slide 21: "... invoke solveProblem on smaller problems ..."
#               -> "... invoke solveProblem() on smaller problems ..."
#   - Slide 21: "... solve all the smaller problem, ..."
#               -> "... solve all the smaller problems, ..."
#
# Edits are made via TextRange.Characters(start, length) on an exact,
# uniquely-located substring so only the affected run(s) are rewritten -
# the rest of the paragraph/run structure (line breaks, indent levels,
# per-run formatting) is left completely untouched.

$p = $ppt.ActivePresentation

function Set-SubstringText($TextRange, $Find, $Replace) {
    $full = $TextRange.Text
    $idx = $full.IndexOf($Find)
    if ($idx -lt 0) {
        throw "Could not find target text: $Find"
    }
    $chars = $TextRange.Characters($idx + 1, $Find.Length)
    $chars.Text = $Replace
}

# --- Slide 12: "An Aside: Proving it right with Loop Invariants" ---
$s12 = $p.Slides.Item(12)
for ($i = 1; $i -le $s12.Shapes.Count; $i++) {
    $shp = $s12.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text.IndexOf("A[0..j-1] are the elements originally stored in the sub-list but in sorted order") -ge 0) {
            Set-SubstringText $tr `
                "A[0..j-1] are the elements originally stored in the sub-list but in sorted order" `
                "A[1..j-1] are the elements originally stored in the sub-list but in sorted order"
        }
    }
}

# --- Slide 21: "Divide and Conquer Strategy" ---
$s21 = $p.Slides.Item(21)
for ($i = 1; $i -le $s21.Shapes.Count; $i++) {
    $shp = $s21.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        $txt = $tr.Text

        if ($txt.IndexOf(" on smaller problems ") -ge 0) {
            Set-SubstringText $tr " on smaller problems " "() on smaller problems "
        }

        # Re-read, since the text range may have shifted after the edit above.
        $txt = $tr.Text
        if ($txt.IndexOf("Note:  maybe solve all the smaller problem, or maybe just some of them.") -ge 0) {
            Set-SubstringText $tr `
                "Note:  maybe solve all the smaller problem, or maybe just some of them." `
                "Note:  maybe solve all the smaller problems, or maybe just some of them."
        }
    }
}
